# Update summary charts / summary reports per Prof. Erhardt's comments.
# Tampa-St. Petersburg-Clearwater, FL Metro Area-Rail FAC summary report.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Window geometry on the workbook (best effort - cosmetic only).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Width  = 29040
$win.Height = 15840
$win.Left   = 28680
$win.Top    = -120

# ---------------------------------------------------------------------------
# 2. Year headers: Year 1 moved from 2002 to 2012.
# ---------------------------------------------------------------------------
$ws1.Range("C1").Value = "2012"
$ws1.Range("E7").Value = "2012"

# ---------------------------------------------------------------------------
# 3. Updated factor table data (Average Values columns E/F and Ridership
#    Effect column H) for rows 8-19.
# ---------------------------------------------------------------------------
$ws1.Range("E8").Value = 67507
$ws1.Range("F8").Value = 75123
$ws1.Range("H8").Value = -5470.054809700003

$ws1.Range("E9").Value = 2.061071585
$ws1.Range("F9").Value = 1.242708385
$ws1.Range("H9").Value = -127045.3087837

$ws1.Range("E10").Value = 4137582.92
$ws1.Range("F10").Value = 4635160.05
$ws1.Range("H10").Value = 14291.654482

$ws1.Range("E11").Value = 16.48231124
$ws1.Range("F11").Value = 16.16132123
$ws1.Range("H11").Value = -527.6458488899999

$ws1.Range("E12").Value = 3.9458
$ws1.Range("F12").Value = 2.71
$ws1.Range("H12").Value = -13871.3142601

$ws1.Range("E13").Value = 27302.32
$ws1.Range("F13").Value = 28301.5
$ws1.Range("H13").Value = -2110.71123519

$ws1.Range("E14").Value = 7.73
$ws1.Range("F14").Value = 6.39
$ws1.Range("H14").Value = -2851.78714502

$ws1.Range("E15").Value = 5.4
$ws1.Range("F15").Value = 7.575
$ws1.Range("H15").Value = -2194.4003394

# Rows 16-18 keep their existing values (E16/F16 blank, E17/F17 = 0/1,
# E18/F18 = 0/0, H16/H17/H18 unchanged) - only formulas/number formats
# change for them below.

# Row 19 "New Reporters" - Riddership Effect now explicitly zero.
$ws1.Range("H19").Value = 0

# Row 20 "Total Modeled Ridership"
$ws1.Range("E20").Value = 320350.1669
$ws1.Range("F20").Value = 362316.8367

# Row 21 "Total Observed Ridership"
$ws1.Range("E21").Value = 301516
$ws1.Range("F21").Value = 423123

# ---------------------------------------------------------------------------
# 4. Formulas: drop the "*100" multiplier from the % Diff / Ridership Effect
#    % formulas now that the cells carry a genuine percentage number format.
# ---------------------------------------------------------------------------
for ($r = 8; $r -le 21; $r++) {
    $ws1.Range("G$r").Formula = "=IFERROR((F$r-E$r)/E$r,0)"
}
foreach ($r in 8..19) {
    $ws1.Range("I$r").Formula = "=IFERROR(H$r/`$E`$21,0)"
}
$ws1.Range("I21").Formula = "=IFERROR(H21/`$E`$21,0)"

# ---------------------------------------------------------------------------
# 5. Number formats: the Average Values (E/F) and Ridership Effect (H)
#    columns now use an accounting-style 2-decimal format, and the % Diff
#    columns (G/I) now use a real percentage format instead of baking the
#    *100 into the formula.
# ---------------------------------------------------------------------------
$ws1.Range("E8:F19").NumberFormat  = "#,##0.00"
$ws1.Range("H8:H19").NumberFormat  = "#,##0.00"
$ws1.Range("G8:G21").NumberFormat  = "0.00%"
$ws1.Range("I8:I21").NumberFormat  = "0.00%"
$ws1.Range("E20:F20").NumberFormat = "#,##0.00"
$ws1.Range("H20").NumberFormat     = "#,##0.00"
$ws1.Range("E21:F21").NumberFormat = "#,##0.00"
$ws1.Range("H21").NumberFormat     = "#,##0.00"

# ---------------------------------------------------------------------------
# 6. Sheet view: scroll position resets and the active selection moves from
#    K20 to H21.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1").Select()
$ws1.Range("H21").Select()
